$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -10
$ws.Range("F3").Value = -7
$ws.Range("F6").Value = 5
$ws.Range("F8").Value = -1
